$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "Förändrad" (column C) timestamp refresh: 45184 -> 45186 for every data row ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# --- 2) Add the friendly display-name second argument to the HYPERLINK() formulas ---
# Map of column letter -> (subfolder, extension) used to build each link.
$linkCols = @(
    @{ Col = "S"; Sub = "artfynd";         Ext = ".xlsx" },
    @{ Col = "T"; Sub = "kartor";          Ext = ".png"  },
    @{ Col = "V"; Sub = "klagomål";        Ext = ".docx" },
    @{ Col = "W"; Sub = "klagomålsmail";   Ext = ".docx" },
    @{ Col = "X"; Sub = "tillsyn";         Ext = ".docx" },
    @{ Col = "Y"; Sub = "tillsynsmail";    Ext = ".docx" }
)

for ($r = 2; $r -le $lastRow; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($beteckning)) {
        continue
    }

    foreach ($lc in $linkCols) {
        $cell = $ws.Range($lc.Col + $r)
        $formula = $cell.Formula
        if ($formula -and $formula -like "*HYPERLINK(*" -and $formula -notlike "*,*") {
            $url = "https://klasma.github.io/Logging_BASTAD/" + $lc.Sub + "/" + $beteckning + $lc.Ext
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
        }
    }
}
